# Applies the ObjectArrange.xlsx edit: populates Sheet1 with the grid-arrangement
# table (ObjectID / Gridx / Gridy per Stage) that the commit introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "ObjectID"
$ws.Range("B1").Value = "Gridx"
$ws.Range("C1").Value = "Gridy"

# Data rows
$ws.Range("A2").Value = 1

$ws.Range("A3").Value = "Stage1"
$ws.Range("B3").Value = 10

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = 2

$ws.Range("A5").Value = 2
$ws.Range("B5").Value = 3
$ws.Range("C5").Value = 3

$ws.Range("A6").Value = 2
$ws.Range("B6").Value = 4
$ws.Range("C6").Value = 4

$ws.Range("A7").Value = 2
$ws.Range("B7").Value = 5
$ws.Range("C7").Value = 5

$ws.Range("A8").Value = 99
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 0

$ws.Range("A9").Value = 3
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = 0

$ws.Range("A10").Value = 3
$ws.Range("B10").Value = 2
$ws.Range("C10").Value = 0

$ws.Range("A11").Value = 4
$ws.Range("B11").Value = 3
$ws.Range("C11").Value = 0

$ws.Range("A12").Value = 4
$ws.Range("B12").Value = 4
$ws.Range("C12").Value = 0

$ws.Range("A13").Value = 5
$ws.Range("B13").Value = 10
$ws.Range("C13").Value = 3

# Page setup: A4, portrait
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Selection should land on the last-edited cell, C13
$ws.Range("C13").Select()
